$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PBL2")

# PBI 12 ("... Forecast bis zu welchem Sprint ...") now forecasts completion
# across two sprints instead of just "Sprint 3".
$ws.Range("F8").Value = "Sprint 3, Sprint 4"

# Move the active selection to F9, matching the sprint-planning review.
$ws.Range("F9").Select()
